$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Cells.Item(6, 9).Value = "aa"
$ws.Cells.Item(6, 10).Value = "Agree/Accept"
$ws.Cells.Item(19, 9).Value = "aa"
$ws.Cells.Item(19, 10).Value = "Agree/Accept"
$ws.Cells.Item(37, 9).Value = "aa"
$ws.Cells.Item(37, 10).Value = "Agree/Accept"
$ws.Cells.Item(48, 9).Value = "sd"
$ws.Cells.Item(48, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(49, 9).Value = "sv"
$ws.Cells.Item(49, 10).Value = "Statement-opinion"
$ws.Cells.Item(56, 9).Value = "aa"
$ws.Cells.Item(56, 10).Value = "Agree/Accept"
$ws.Cells.Item(68, 9).Value = "b"
$ws.Cells.Item(68, 10).Value = "Acknowledge (Backchannel)"
$ws.Cells.Item(71, 9).Value = "sd"
$ws.Cells.Item(71, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(79, 9).Value = "sv"
$ws.Cells.Item(79, 10).Value = "Statement-opinion"
$ws.Cells.Item(84, 9).Value = "aa"
$ws.Cells.Item(84, 10).Value = "Agree/Accept"
$ws.Cells.Item(85, 9).Value = "aa"
$ws.Cells.Item(85, 10).Value = "Agree/Accept"
$ws.Cells.Item(100, 9).Value = "sv"
$ws.Cells.Item(100, 10).Value = "Statement-opinion"
$ws.Cells.Item(104, 9).Value = "sv"
$ws.Cells.Item(104, 10).Value = "Statement-opinion"
$ws.Cells.Item(112, 9).Value = "sd"
$ws.Cells.Item(112, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(113, 9).Value = "sd"
$ws.Cells.Item(113, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(126, 9).Value = "sd"
$ws.Cells.Item(126, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(138, 9).Value = "qy"
$ws.Cells.Item(138, 10).Value = "Yes-No-Question"
$ws.Cells.Item(141, 9).Value = "%"
$ws.Cells.Item(141, 10).Value = "Uninterpretable"
$ws.Cells.Item(164, 9).Value = "%"
$ws.Cells.Item(164, 10).Value = "Uninterpretable"
$ws.Cells.Item(186, 9).Value = "aa"
$ws.Cells.Item(186, 10).Value = "Agree/Accept"
$ws.Cells.Item(187, 9).Value = "aa"
$ws.Cells.Item(187, 10).Value = "Agree/Accept"
$ws.Cells.Item(190, 9).Value = "aa"
$ws.Cells.Item(190, 10).Value = "Agree/Accept"
$ws.Cells.Item(209, 9).Value = "sd"
$ws.Cells.Item(209, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(214, 9).Value = "aa"
$ws.Cells.Item(214, 10).Value = "Agree/Accept"
$ws.Cells.Item(236, 9).Value = "aa"
$ws.Cells.Item(236, 10).Value = "Agree/Accept"
$ws.Cells.Item(238, 9).Value = "sd"
$ws.Cells.Item(238, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(246, 9).Value = "%"
$ws.Cells.Item(246, 10).Value = "Uninterpretable"
$ws.Cells.Item(262, 9).Value = "aa"
$ws.Cells.Item(262, 10).Value = "Agree/Accept"
$ws.Cells.Item(265, 9).Value = "sd"
$ws.Cells.Item(265, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(268, 9).Value = "sd"
$ws.Cells.Item(268, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(273, 9).Value = "sd"
$ws.Cells.Item(273, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(296, 9).Value = "sd"
$ws.Cells.Item(296, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(299, 9).Value = "sd"
$ws.Cells.Item(299, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(304, 9).Value = "aa"
$ws.Cells.Item(304, 10).Value = "Agree/Accept"
$ws.Cells.Item(307, 9).Value = "sv"
$ws.Cells.Item(307, 10).Value = "Statement-opinion"
$ws.Cells.Item(318, 9).Value = "sd"
$ws.Cells.Item(318, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(323, 9).Value = "sd"
$ws.Cells.Item(323, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(328, 9).Value = "sd"
$ws.Cells.Item(328, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(337, 9).Value = "aa"
$ws.Cells.Item(337, 10).Value = "Agree/Accept"
$ws.Cells.Item(339, 9).Value = "aa"
$ws.Cells.Item(339, 10).Value = "Agree/Accept"
$ws.Cells.Item(343, 9).Value = "%"
$ws.Cells.Item(343, 10).Value = "Uninterpretable"
